$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = 1086.29
$ws.Range("O4").Value = 778.14
$ws.Range("O5").Value = 1682.34
$ws.Range("O8").Value = 996.08
$ws.Range("O11").Value = 1383.28
$ws.Range("O12").Value = 1578.51
$ws.Range("O13").Value = 1218.36
$ws.Range("O14").Value = 1183.95
$ws.Range("O15").Value = 1086.8
$ws.Range("O18").Value = 1245.16
$ws.Range("O20").Value = 1016.08
$ws.Range("O21").Value = 1045.29
$ws.Range("O23").Value = 1314.03
$ws.Range("O24").Value = 1082.64
$ws.Range("O26").Value = 1133.68
$ws.Range("O28").Value = 902.71
$ws.Range("O30").Value = 1354.34
$ws.Range("O34").Value = 1283.92
$ws.Range("O35").Value = 1008.44
$ws.Range("O37").Value = 894.88
$ws.Range("O39").Value = 1070.6
$ws.Range("O41").Value = 1509.14
$ws.Range("O42").Value = 1225.43
$ws.Range("O44").Value = 1026.07
$ws.Range("O45").Value = 1058.49
$ws.Range("O46").Value = 1059.44
$ws.Range("O49").Value = 1238.38
$ws.Range("O51").Value = 878.3200000000001
$ws.Range("O54").Value = 1069.91
$ws.Range("O57").Value = 1230.59
$ws.Range("O59").Value = 902.91
$ws.Range("O60").Value = 813.87
$ws.Range("O61").Value = 1636.25
$ws.Range("O63").Value = 1027.84
$ws.Range("O66").Value = 1789.55
$ws.Range("O67").Value = 1039.42
$ws.Range("O68").Value = 1064.4
$ws.Range("O69").Value = 1296.26
$ws.Range("O72").Value = 1114.56
$ws.Range("O74").Value = 1085.29
$ws.Range("O79").Value = 1254.98
$ws.Range("O82").Value = 1108.31
$ws.Range("O84").Value = 1526.67
$ws.Range("O85").Value = 1467.07
$ws.Range("O87").Value = 1179.97
$ws.Range("O88").Value = 1452.54
$ws.Range("O89").Value = 1361.81
$ws.Range("O90").Value = 1415.12
$ws.Range("O91").Value = 667.23
$ws.Range("O92").Value = 1368.65
$ws.Range("O95").Value = 1430.34
$ws.Range("O96").Value = 1059.23
$ws.Range("O97").Value = 1524.09
$ws.Range("O99").Value = 825.0700000000001
$ws.Range("O100").Value = 1405.18
$ws.Range("O102").Value = 1163.77
$ws.Range("O103").Value = 1363.56
$ws.Range("O104").Value = 886.7
$ws.Range("O106").Value = 1291.04
$ws.Range("O108").Value = 1306.63
$ws.Range("O109").Value = 1361.13
$ws.Range("O110").Value = 1781.13
$ws.Range("O111").Value = 1285.9
$ws.Range("O112").Value = 977.23
$ws.Range("O113").Value = 1095.58
$ws.Range("O114").Value = 1521.41
$ws.Range("O115").Value = 1152.08
$ws.Range("O117").Value = 1565.96
$ws.Range("O124").Value = 1102.93
$ws.Range("O125").Value = 1124.82
$ws.Range("O126").Value = 833.17
$ws.Range("O127").Value = 1281.33
$ws.Range("O129").Value = 972.6900000000001
$ws.Range("O130").Value = 1390.83
$ws.Range("O133").Value = 1282.58
$ws.Range("O136").Value = 1521.97
$ws.Range("O139").Value = 1121.08
$ws.Range("O143").Value = 1098.1
$ws.Range("O145").Value = 1425.61
$ws.Range("O148").Value = 1447.38
$ws.Range("O149").Value = 1234.43
$ws.Range("O150").Value = 1673.21
$ws.Range("O151").Value = 1490.79
$ws.Range("O152").Value = 1386.66
$ws.Range("O153").Value = 1238.23
$ws.Range("O154").Value = 1637.9
$ws.Range("O155").Value = 1429.94
$ws.Range("O156").Value = 1567.58
$ws.Range("O157").Value = 1081.06
$ws.Range("O161").Value = 1263.11
$ws.Range("O163").Value = 1141.47
$ws.Range("O165").Value = 1183.28
$ws.Range("O170").Value = 1050.03
$ws.Range("O171").Value = 1394.98
$ws.Range("O173").Value = 1156.34
$ws.Range("O174").Value = 1299.4
$ws.Range("O178").Value = 1113.97
$ws.Range("O184").Value = 1415.99
$ws.Range("O185").Value = 1009.05
$ws.Range("O187").Value = 1074.15
$ws.Range("O190").Value = 1213.65
$ws.Range("O191").Value = 909.96
$ws.Range("O192").Value = 1399.5
$ws.Range("O196").Value = 1645.95
$ws.Range("O197").Value = 1322.37
$ws.Range("O199").Value = 1222.16
$ws.Range("O201").Value = 1200.16
$ws.Range("O202").Value = 1388.54
$ws.Range("O203").Value = 893.04
$ws.Range("O204").Value = 1315.84
$ws.Range("O205").Value = 1412.61
$ws.Range("O207").Value = 1483.71
$ws.Range("O208").Value = 1499.77
$ws.Range("O209").Value = 957.73
$ws.Range("O210").Value = 1220.46
$ws.Range("O211").Value = 1368.23
$ws.Range("O212").Value = 1139.17
$ws.Range("O213").Value = 1023.38
$ws.Range("O215").Value = 1134.27
$ws.Range("O217").Value = 1362.71
$ws.Range("O219").Value = 1088.28
$ws.Range("O220").Value = 1718.08
$ws.Range("O221").Value = 1105.35
$ws.Range("O223").Value = 1345.73
$ws.Range("O226").Value = 1141.99
$ws.Range("O227").Value = 860.4400000000001
$ws.Range("O229").Value = 1288.79
$ws.Range("O231").Value = 1223.92
